# Actualización automática del tracker
# Rellena las columnas G (resultado) y H (profit) para las filas recién
# resueltas, dejando intactas las filas aún sin resultado (65 y 66).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 61; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 62; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 63; Resultado = "Acierto"; Profit = 1.5 },
    @{ Row = 64; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 67; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 68; Resultado = "Acierto"; Profit = 0.91 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Resultado
    $ws.Cells.Item($u.Row, 8).Value = $u.Profit
}

$wb.Save()
